# Add a "min_units" column (C) to the requirements/course table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1 - same label + style (bold/border/centered) as A1/B1.
$ws.Range("C1").Value = "min_units"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row minimum-units values (rows 2..96), aligned with existing A/B data.
$minUnits = @(10,10,10,10,10,10,12,10,10,1,9,9,12,12,12,12,12,12,12,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,12,12,12,12,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9)

for ($i = 0; $i -lt $minUnits.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $minUnits[$i]
}

Write-Host "min_units column written"
